$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 48 (A48="-544", "Vera 453") is removed; subsequent rows shift up by one.
$ws.Rows.Item(48).Delete()
